# Lab3 Timeseries powerflow - vm_pu.xlsx results update
# Rewrites the voltage-magnitude (p.u.) result table for rows 2-25 (timesteps 0-23)
# with the completed run's values, and fills in the previously-missing
# bus columns I, J, K (buses 7, 8, 9) that the earlier (partial) run had skipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.991412918452992
$ws.Range("D2").Value = 0.986601626550307
$ws.Range("E2").Value = 0.992613680192833
$ws.Range("F2").Value = 0.9999999999999999
$ws.Range("G2").Value = 0.974316325861903
$ws.Range("H2").Value = 0.9874975009839858
$ws.Range("I2").Value = 1.00191104615943
$ws.Range("J2").Value = 0.9868976472291071
$ws.Range("K2").Value = 0.997709920870801
$ws.Range("L2").Value = 0.991071574230921
$ws.Range("M2").Value = 0.9854014881851936
$ws.Range("N2").Value = 0.9958068879420655
$ws.Range("O2").Value = 0.9706684610392248
$ws.Range("C3").Value = 0.9927068647806099
$ws.Range("D3").Value = 0.9884466247576442
$ws.Range("E3").Value = 0.9938326202608346
$ws.Range("G3").Value = 0.9771696028160756
$ws.Range("H3").Value = 0.9893441745247911
$ws.Range("I3").Value = 1.00191104615943
$ws.Range("J3").Value = 0.9887431990110654
$ws.Range("K3").Value = 0.998983502820711
$ws.Range("L3").Value = 0.992413460816302
$ws.Range("M3").Value = 0.9873358439731758
$ws.Range("N3").Value = 0.9970780406604551
$ws.Range("O3").Value = 0.9738434027420182
$ws.Range("C4").Value = 0.9876275637982947
$ws.Range("D4").Value = 0.9806856990150333
$ws.Range("E4").Value = 0.9889005695763177
$ws.Range("G4").Value = 0.964521919088422
$ws.Range("H4").Value = 0.98157620154572
$ws.Range("I4").Value = 1.00191104615943
$ws.Range("J4").Value = 0.9809799446745778
$ws.Range("K4").Value = 0.9940227554399166
$ws.Range("L4").Value = 0.9873324329905674
$ws.Range("M4").Value = 0.9791433665504153
$ws.Range("N4").Value = 0.9921267554143143
$ws.Range("O4").Value = 0.9595980535122357
$ws.Range("C5").Value = 0.9880405230101313
$ws.Range("D5").Value = 0.9813115129876749
$ws.Range("E5").Value = 0.9892982900713927
$ws.Range("G5").Value = 0.9655416561908783
$ws.Range("H5").Value = 0.9822025837829206
$ws.Range("I5").Value = 1.00191104615943
$ws.Range("J5").Value = 0.981605946416905
$ws.Range("K5").Value = 0.9944227909938281
$ws.Range("L5").Value = 0.987745533626158
$ws.Range("M5").Value = 0.9798037843388554
$ws.Range("N5").Value = 0.9925260279399992
$ws.Range("O5").Value = 0.9607499660356039
$ws.Range("C6").Value = 0.9825674952481146
$ws.Range("D6").Value = 0.9759489442219574
$ws.Range("E6").Value = 0.9850001995921686
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9599594632248213
$ws.Range("H6").Value = 0.9768351455864962
$ws.Range("I6").Value = 1.00191104615943
$ws.Range("J6").Value = 0.9762417686621047
$ws.Range("K6").Value = 0.9890080593965691
$ws.Range("L6").Value = 0.9811754057837218
$ws.Range("M6").Value = 0.974422529401186
$ws.Range("N6").Value = 0.9871216244073553
$ws.Range("O6").Value = 0.9550916407657314
$ws.Range("C7").Value = 0.991011722163194
$ws.Range("D7").Value = 0.9902495639652201
$ws.Range("E7").Value = 0.9933027343024717
$ws.Range("F7").Value = 0.9999999999999999
$ws.Range("G7").Value = 0.984741140893128
$ws.Range("H7").Value = 0.9911487508745511
$ws.Range("I7").Value = 1.00191104615943
$ws.Range("J7").Value = 0.9905466791738494
$ws.Range("K7").Value = 0.996177732192194
$ws.Range("L7").Value = 0.9884359867814143
$ws.Range("M7").Value = 0.9896517055774779
$ws.Range("N7").Value = 0.9942776217617193
$ws.Range("O7").Value = 0.9831481071063324
$ws.Range("C8").Value = 0.9883965474584619
$ws.Range("D8").Value = 0.9870300066223741
$ws.Range("E8").Value = 0.9909931175490267
$ws.Range("F8").Value = 0.9999999999999999
$ws.Range("G8").Value = 0.9803816681823642
$ws.Range("H8").Value = 0.9879262700425949
$ws.Range("I8").Value = 1.00191104615943
$ws.Range("J8").Value = 0.9873261558326467
$ws.Range("K8").Value = 0.9934294259627219
$ws.Range("L8").Value = 0.9853903811192537
$ws.Range("M8").Value = 0.9863312710194294
$ws.Range("N8").Value = 0.9915345576543744
$ws.Range("O8").Value = 0.9784630782537093
$ws.Range("C9").Value = 0.9897890356080422
$ws.Range("D9").Value = 0.9884632155952361
$ws.Range("E9").Value = 0.9921646160006153
$ws.Range("G9").Value = 0.9819496678335841
$ws.Range("H9").Value = 0.9893607804275392
$ws.Range("I9").Value = 1.00191104615943
$ws.Range("J9").Value = 0.9887597948265848
$ws.Range("K9").Value = 0.9950272184133685
$ws.Range("L9").Value = 0.9872066241273124
$ws.Range("M9").Value = 0.9877763053471029
$ws.Range("N9").Value = 0.9931293024740582
$ws.Range("O9").Value = 0.9800701098586213
$ws.Range("C10").Value = 0.9826119284642001
$ws.Range("D10").Value = 0.9771630686930983
$ws.Range("E10").Value = 0.9853521020509547
$ws.Range("F10").Value = 0.9999999999999999
$ws.Range("G10").Value = 0.9631908147892677
$ws.Range("H10").Value = 0.9780503725320742
$ws.Range("I10").Value = 1.00191104615943
$ws.Range("J10").Value = 0.9774562574200464
$ws.Range("K10").Value = 0.9887492927963573
$ws.Range("L10").Value = 0.9806029195305277
$ws.Range("M10").Value = 0.9758145109284186
$ws.Range("N10").Value = 0.9868633513788228
$ws.Range("O10").Value = 0.9589924310616695
$ws.Range("C11").Value = 0.9801251315318334
$ws.Range("D11").Value = 0.9737174816890366
$ws.Range("E11").Value = 0.9831088128833123
$ws.Range("G11").Value = 0.9579094060441602
$ws.Range("H11").Value = 0.9746016567947702
$ws.Range("I11").Value = 1.00191104615943
$ws.Range("J11").Value = 0.9740096365995217
$ws.Range("K11").Value = 0.9864121498796796
$ws.Range("L11").Value = 0.9780329713119801
$ws.Range("M11").Value = 0.9722073077348008
$ws.Range("N11").Value = 0.9845306663309468
$ws.Range("O11").Value = 0.9531012161532055
$ws.Range("C12").Value = 0.9866011721153686
$ws.Range("D12").Value = 0.9837971484161813
$ws.Range("E12").Value = 0.9892153305456723
$ws.Range("G12").Value = 0.9746085993485624
$ws.Range("H12").Value = 0.9846904762695666
$ws.Range("I12").Value = 1.00191104615943
$ws.Range("J12").Value = 0.9840923276372964
$ws.Range("K12").Value = 0.9920856098494995
$ws.Range("L12").Value = 0.9840413789166247
$ws.Range("M12").Value = 0.9828726883309427
$ws.Range("N12").Value = 0.9901933047373881
$ws.Range("O12").Value = 0.9719317158919567
$ws.Range("C13").Value = 0.9914852324653254
$ws.Range("D13").Value = 0.9912648626430132
$ws.Range("E13").Value = 0.9937859587676097
$ws.Range("F13").Value = 0.9999999999999997
$ws.Range("G13").Value = 0.986674038432259
$ws.Range("H13").Value = 0.9921649714848685
$ws.Range("I13").Value = 1.00191104615943
$ws.Range("J13").Value = 0.9915622824826057
$ws.Range("K13").Value = 0.9963984559869955
$ws.Range("L13").Value = 0.9886428515148532
$ws.Range("M13").Value = 0.9907481576604448
$ws.Range("N13").Value = 0.9944979245477278
$ws.Range("O13").Value = 0.9853351349085974
$ws.Range("C14").Value = 0.9874433406938045
$ws.Range("D14").Value = 0.9872495954449433
$ws.Range("E14").Value = 0.9903533457556403
$ws.Range("G14").Value = 0.9825303383729179
$ws.Range("H14").Value = 0.9881460582607529
$ws.Range("I14").Value = 1.00191104615943
$ws.Range("J14").Value = 0.9875458105408079
$ws.Range("K14").Value = 0.9914549780550801
$ws.Range("L14").Value = 0.9830938245858712
$ws.Range("M14").Value = 0.9867219790690366
$ws.Range("N14").Value = 0.9895638758107015
$ws.Range("O14").Value = 0.9811566432922906
$ws.Range("C15").Value = 0.9923527992991032
$ws.Range("D15").Value = 0.992361309056801
$ws.Range("E15").Value = 0.9945583359714845
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.9881819676922552
$ws.Range("H15").Value = 0.9932624135166275
$ws.Range("I15").Value = 1.00191104615943
$ws.Range("J15").Value = 0.9926590578749834
$ws.Range("K15").Value = 0.9972946518887379
$ws.Range("L15").Value = 0.9896312332750119
$ws.Range("M15").Value = 0.9918808349578688
$ws.Range("N15").Value = 0.9953924110444857
$ws.Range("O15").Value = 0.9869543838822524
$ws.Range("C16").Value = 0.9963332161159905
$ws.Range("D16").Value = 0.9972589630574911
$ws.Range("E16").Value = 0.9981675300023829
$ws.Range("G16").Value = 0.9947178922364165
$ws.Range("H16").Value = 0.9981645147864948
$ws.Range("I16").Value = 1.00191104615943
$ws.Range("J16").Value = 0.9975581813713874
$ws.Range("K16").Value = 1.001699594839977
$ws.Range("L16").Value = 0.9943929646026743
$ws.Range("M16").Value = 0.9969218825945572
$ws.Range("N16").Value = 0.999788952003011
$ws.Range("O16").Value = 0.9939146661144467
$ws.Range("C17").Value = 0.9974500724858488
$ws.Range("D17").Value = 0.9987756254954264
$ws.Range("E17").Value = 0.9992287038255466
$ws.Range("G17").Value = 0.9969069038750491
$ws.Range("H17").Value = 0.9996825544156555
$ws.Range("I17").Value = 1.00191104615943
$ws.Range("J17").Value = 0.9990752988698381
$ws.Range("K17").Value = 1.002912797241131
$ws.Range("L17").Value = 0.9956567791675009
$ws.Range("M17").Value = 0.9984969518749355
$ws.Range("N17").Value = 1.000999840340659
$ws.Range("O17").Value = 0.9962677088875616
$ws.Range("C18").Value = 0.9845670539985825
$ws.Range("D18").Value = 0.9798106712454032
$ws.Range("E18").Value = 0.987093069210654
$ws.Range("G18").Value = 0.9671536092193526
$ws.Range("H18").Value = 0.9807003792152598
$ws.Range("I18").Value = 1.00191104615943
$ws.Range("J18").Value = 0.9801046543609714
$ws.Range("K18").Value = 0.9905843836833218
$ws.Range("L18").Value = 0.9826440001248464
$ws.Range("M18").Value = 0.978578306545616
$ws.Range("N18").Value = 0.9886949420115431
$ws.Range("O18").Value = 0.9633840497792517
$ws.Range("C19").Value = 0.9492824053062057
$ws.Range("D19").Value = 0.9340172462403051
$ws.Range("E19").Value = 0.9565060656078134
$ws.Range("F19").Value = 0.9999999999999999
$ws.Range("G19").Value = 0.8998564340281798
$ws.Range("H19").Value = 0.9348653719163675
$ws.Range("I19").Value = 1.00191104615943
$ws.Range("J19").Value = 0.9342974894629013
$ws.Range("K19").Value = 0.957426355164789
$ws.Range("L19").Value = 0.9448901694080496
$ws.Range("M19").Value = 0.9309563505054248
$ws.Range("N19").Value = 0.9556001591506935
$ws.Range("O19").Value = 0.8884087971500502
$ws.Range("C20").Value = 0.917597626130035
$ws.Range("D20").Value = 0.9027194897537949
$ws.Range("E20").Value = 0.9320689858744271
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.8669606608930009
$ws.Range("H20").Value = 0.9035391957930828
$ws.Range("I20").Value = 1.00191104615943
$ws.Range("J20").Value = 0.9029903423744799
$ws.Range("K20").Value = 0.9241201716847359
$ws.Range("L20").Value = 0.9039826685193088
$ws.Range("M20").Value = 0.8995403708963997
$ws.Range("N20").Value = 0.9223575039192498
$ws.Range("O20").Value = 0.8548475557083314
$ws.Range("C21").Value = 0.9293500909875693
$ws.Range("D21").Value = 0.9138804091935714
$ws.Range("E21").Value = 0.9410640653957161
$ws.Range("F21").Value = 0.9999999999999999
$ws.Range("G21").Value = 0.8778244172367532
$ws.Range("H21").Value = 0.9147102498019837
$ws.Range("I21").Value = 1.00191104615943
$ws.Range("J21").Value = 0.9141546105447466
$ws.Range("K21").Value = 0.9369952426816992
$ws.Range("L21").Value = 0.9198270535410128
$ws.Range("M21").Value = 0.9106735667638671
$ws.Range("N21").Value = 0.9352080169925575
$ws.Range("O21").Value = 0.8656173274125788
$ws.Range("C22").Value = 0.9407719281693283
$ws.Range("D22").Value = 0.9226243133512289
$ws.Range("E22").Value = 0.9490935098662225
$ws.Range("G22").Value = 0.8823937265654027
$ws.Range("H22").Value = 0.9234620937805118
$ws.Range("I22").Value = 1.00191104615943
$ws.Range("J22").Value = 0.9229011382298486
$ws.Range("K22").Value = 0.9496949750233646
$ws.Range("L22").Value = 0.9360823254740485
$ws.Range("M22").Value = 0.9190840058107757
$ws.Range("N22").Value = 0.9478835258516989
$ws.Range("O22").Value = 0.8685842642186784
$ws.Range("C23").Value = 0.9509914511532186
$ws.Range("D23").Value = 0.9331867504943965
$ws.Range("E23").Value = 0.957000582092611
$ws.Range("F23").Value = 0.9999999999999999
$ws.Range("G23").Value = 0.8944736197960667
$ws.Range("H23").Value = 0.9340341220465187
$ws.Range("I23").Value = 1.00191104615943
$ws.Range("J23").Value = 0.9334667445344373
$ws.Range("K23").Value = 0.9596381406166152
$ws.Range("L23").Value = 0.9483545814827352
$ws.Range("M23").Value = 0.9297614513280639
$ws.Range("N23").Value = 0.9578077258406753
$ws.Range("O23").Value = 0.8812764698440747
$ws.Range("C24").Value = 0.9559145884041806
$ws.Range("D24").Value = 0.9369990347425817
$ws.Range("E24").Value = 0.9603141152047258
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.8966637194674282
$ws.Range("H24").Value = 0.9378498680039694
$ws.Range("I24").Value = 1.00191104615943
$ws.Range("J24").Value = 0.9372801726231968
$ws.Range("K24").Value = 0.9642879752634341
$ws.Range("L24").Value = 0.9546190979063724
$ws.Range("M24").Value = 0.9334444465651944
$ws.Range("N24").Value = 0.9624486913881086
$ws.Range("O24").Value = 0.8828430588732281
$ws.Range("C25").Value = 0.9627094657934155
$ws.Range("D25").Value = 0.9452807128533109
$ws.Range("E25").Value = 0.96599051929267
$ws.Range("F25").Value = 0.9999999999999998
$ws.Range("G25").Value = 0.9083447246465038
$ws.Range("H25").Value = 0.9461390662155049
$ws.Range("I25").Value = 1.00191104615943
$ws.Range("J25").Value = 0.945564335574727
$ws.Range("K25").Value = 0.970527862011329
$ws.Range("L25").Value = 0.961952397139046
$ws.Range("M25").Value = 0.9419924006876205
$ws.Range("N25").Value = 0.9686766761696058
$ws.Range("O25").Value = 0.8958561496377084
